$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell reference -> new value, derived from the authoritative diff
# (prices in column D, percentage volumes in column E, plus a swap of the
#  PancakeSwap / EthereumClassic rows 32-33 including their B/C/D/E data).
$updates = [ordered]@{
    'D2' = '59.063.30'
    'E2' = '  +3.92%  '
    'D3' = '2.595.77'
    'E3' = '  +2.31%  '
    'E4' = '  +0.10%  '
    'D5' = '521.56'
    'E5' = '  +1.77%  '
    'D6' = '141.37'
    'E6' = '  +1.29%  '
    'E7' = '  -0.31%  '
    'E8' = '  +2.44%  '
    'D9' = '2.618.92'
    'E9' = '  +3.14%  '
    'E10' = '  +1.06%  '
    'E11' = '  +2.36%  '
    'E12' = '  +3.04%  '
    'E13' = '  +2.53%  '
    'D14' = '3.058.69'
    'E14' = '  +2.59%  '
    'D15' = '59.125.43'
    'E15' = '  +3.94%  '
    'D16' = '20.58'
    'E16' = '  +2.88%  '
    'D17' = '2.620.37'
    'E17' = '  +3.89%  '
    'E18' = '  +0.72%  '
    'D19' = '339.77'
    'E19' = '  +2.60%  '
    'E20' = '  +1.70%  '
    'D21' = '10.23'
    'E21' = '  +1.77%  '
    'D22' = '6.58'
    'E22' = '  +7.58%  '
    'D23' = '0.997'
    'E23' = '  -0.28%  '
    'D24' = '66.37'
    'E24' = '  +3.66%  '
    'E25' = '  +1.99%  '
    'E26' = '  +1.84%  '
    'E27' = '  -0.35%  '
    'E28' = '  +3.99%  '
    'D29' = '0.998'
    'E29' = '  -0.11%  '
    'D30' = '0.0₃0729'
    'E30' = '  -2.49%  '
    'D31' = '5.97'
    'E31' = '  -4.19%  '
    'B32' = 'PancakeSwap'
    'C32' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D32' = '1.57'
    'E32' = '  +1.60%  '
    'B33' = 'EthereumClassic'
    'C33' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D33' = '18.80'
    'E33' = '  +2.17%  '
    'D34' = '149.00'
    'E34' = '  +0.65%  '
    'D35' = '4.02'
    'E35' = '  +1.49%  '
    'E36' = '  +1.06%  '
    'D37' = '36.34'
    'E37' = '  +1.99%  '
    'E38' = '  +3.63%  '
    'D39' = '0.837'
    'E39' = '  +2.04%  '
    'D40' = '0.831'
    'E40' = '  -0.92%  '
    'E41' = '  +3.07%  '
    'D42' = '277.50'
    'E42' = '  +7.02%  '
    'E43' = '  -0.49%  '
    'E44' = '  +1.18%  '
    'E45' = '  +0.60%  '
    'E46' = '  +2.75%  '
    'E47' = '  +1.06%  '
    'D48' = '18.67'
    'E48' = '  +1.40%  '
    'D49' = '1.988.76'
    'E49' = '  +1.30%  '
    'D50' = '4.63'
    'E50' = '  +3.09%  '
    'E51' = '  +0.69%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so Excel does not coerce numeric-looking strings
    # (e.g. '141.37', '18.80', '0.998') into actual numbers, which would drop
    # significant trailing zeros / reformat the text representation.
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
    $cell.Style = 'Normal'
}

Write-Output "Updated $($updates.Count) cells"
